# Update simple accessioning forms
#
# Inserts a new "Tags:" column (with its sample value) between the existing
# "Sequencing Date:" (N) / "Files:" (O->P) columns, i.e. a brand-new column
# is inserted at P, pushing the former P:T columns to Q:U.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at P; everything from P..T shifts right to Q..U.
$ws.Columns("P").Insert() | Out-Null

# Populate the new column's header + first row's sample value.
$ws.Range("P1").Value = "Tags:"
$ws.Range("P2").Value = "t1dm_control, xyz_clinical_trial"

# Match the target column width (≈26.5 "characters" once Excel's standard
# column padding of 5/6 char is added back on save).
$ws.Columns("P").ColumnWidth = 25.666666666666668

# Rows 3 & 4 have nothing in the new Tags column, so make sure no stray
# formatted-but-empty cells are left behind there (Insert() carries over
# the left neighbor's style onto the blank cells it creates).
$ws.Range("P3").ClearContents() | Out-Null
$ws.Range("P4").ClearContents() | Out-Null
$ws.Range("P3").Clear() | Out-Null
$ws.Range("P4").Clear() | Out-Null

# Match the saved selection state.
$ws.Range("P3").Select() | Out-Null
